# Update epexspot_prices.xlsx:
#  - "Prix Spot": a new daily column "15-dec" is inserted right before the
#    "01-oct." column (column ER), shifting all following day columns one
#    column to the right (ER:FV -> ES:FW). The new column has no data yet,
#    so every hourly row (2-25) gets a "-" placeholder, matching the other
#    not-yet-available day columns.
#  - "Gaz" and "CO2": two new daily rows are appended (2025-12-13 and
#    2025-12-14).

$wb = $excel.ActiveWorkbook

# --- "Prix Spot": insert new day column "15-dec" before column ER (148) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns.Item(148).Insert()
$ws1.Range("ER1").Value = "15-dec"
$ws1.Range("ER2:ER25").Value = "-"

# --- "Gaz": append two new daily rows ---
$ws2 = $wb.Worksheets.Item("Gaz")

$ws2.Range("A178").NumberFormat = "@"
$ws2.Range("A178").Value = "2025-12-13"
$ws2.Range("A178").ClearFormats()
$ws2.Range("B178").Value = 26.075

$ws2.Range("A179").NumberFormat = "@"
$ws2.Range("A179").Value = "2025-12-14"
$ws2.Range("A179").ClearFormats()
$ws2.Range("B179").Value = 26.075

# --- "CO2": append two new daily rows ---
$ws3 = $wb.Worksheets.Item("CO2")

$ws3.Range("A178").NumberFormat = "@"
$ws3.Range("A178").Value = "2025-12-13"
$ws3.Range("A178").ClearFormats()
$ws3.Range("B178").Value = 84.09999999999999

$ws3.Range("A179").NumberFormat = "@"
$ws3.Range("A179").Value = "2025-12-14"
$ws3.Range("A179").ClearFormats()
$ws3.Range("B179").Value = 84.09999999999999
